$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.885.24"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.590.12"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "210.01"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.484"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "18.23"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.811.04"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").Value = "1.590.36"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("D16").Value = "25.891.21"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "60.20"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "193.28"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "141.76"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "6.47"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").Value = "1.109.06"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.507"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.34"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "0.783"
$ws.Range("E41").Value = "  -6.41%  "
$ws.Range("D42").Value = "0.822"
$ws.Range("E42").Value = "  +9.87%  "
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("D44").Value = "93.80"
$ws.Range("E44").Value = "  -4.42%  "
$ws.Range("D45").Value = "1.724.53"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "53.56"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  -0.09%  "
